# Insert two new price records (Primera / Segunda) at the top of the
# "Frutilla" Terminal La Palmera de La Serena block, pushing the existing
# rows 791-851 down to 793-853 (dimension grows from A1:T851 to A1:T853).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 791:851 down by two rows.
$ws.Rows("791:792").Insert()

# New row 791 - "Primera" quality entry.
$ws.Cells.Item(791, 1).Value = 8
$ws.Cells.Item(791, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(791, 3).Value = "Coquimbo"
$ws.Cells.Item(791, 4).Value = 44826
$ws.Cells.Item(791, 5).Value = 4
$ws.Cells.Item(791, 6).Value = "Fruta"
$ws.Cells.Item(791, 7).Value = 100101
$ws.Cells.Item(791, 8).Value = "Berries"
$ws.Cells.Item(791, 9).Value = 100112025
$ws.Cells.Item(791, 10).Value = "Frutilla"
$ws.Cells.Item(791, 11).Value = "Sin especificar"
$ws.Cells.Item(791, 12).Value = "Primera"
$ws.Cells.Item(791, 13).Value = 500
$ws.Cells.Item(791, 14).Value = 19000
$ws.Cells.Item(791, 15).Value = 20000
$ws.Cells.Item(791, 16).Value = 19500
$ws.Cells.Item(791, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(791, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(791, 19).Value = 2786
$ws.Cells.Item(791, 20).Value = 7

# New row 792 - "Segunda" quality entry.
$ws.Cells.Item(792, 1).Value = 8
$ws.Cells.Item(792, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(792, 3).Value = "Coquimbo"
$ws.Cells.Item(792, 4).Value = 44826
$ws.Cells.Item(792, 5).Value = 4
$ws.Cells.Item(792, 6).Value = "Fruta"
$ws.Cells.Item(792, 7).Value = 100101
$ws.Cells.Item(792, 8).Value = "Berries"
$ws.Cells.Item(792, 9).Value = 100112025
$ws.Cells.Item(792, 10).Value = "Frutilla"
$ws.Cells.Item(792, 11).Value = "Sin especificar"
$ws.Cells.Item(792, 12).Value = "Segunda"
$ws.Cells.Item(792, 13).Value = 280
$ws.Cells.Item(792, 14).Value = 14000
$ws.Cells.Item(792, 15).Value = 15000
$ws.Cells.Item(792, 16).Value = 14500
$ws.Cells.Item(792, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(792, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(792, 19).Value = 2071
$ws.Cells.Item(792, 20).Value = 7
